$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ClickThroughRateSheet")

# --- Refresh the sample data in rows 3-7 (new Clicks / Impressions numbers) ---
$ws.Range("B3").Value = 423
$ws.Range("C3").Value = 398503

$ws.Range("B4").Value = 757
$ws.Range("C4").Value = 460983

$ws.Range("B5").Value = 235
$ws.Range("C5").Value = 342456

$ws.Range("B6").Value = 145
$ws.Range("C6").Value = 235098

$ws.Range("B7").Value = 134
$ws.Range("C7").Value = 53843

# --- Row 8 keeps only its (now empty) Impressions/CTR cells ---
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()

# --- Row 9 keeps only its (now empty) Impressions cell ---
$ws.Range("B9").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()

# --- Drop the old extra sample rows 10-13 completely ---
$ws.Range("A10:D13").EntireRow.Delete()

# --- Re-establish the shared formula over the remaining data range D3:D7 ---
$ws.Range("D3:D7").Formula = "=B3/C3"
